# inscription BO queue consumer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "BO" as a consumer (column D) for the existing queues it now consumes.
# (Row numbers below are in the ORIGINAL layout, before the new row is inserted.)
$ws.Range("D8").Value = "BO"
$ws.Range("D13").Value = "BO"
$ws.Range("D14").Value = "BO"
$ws.Range("D19").Value = "BO"
$ws.Range("D21").Value = "BO"
$ws.Range("D37").Value = "BO"
$ws.Range("D38").Value = "BO"

# Insert the new CRM_client_<paris> queue row (after "CRM_client") and mark BO as its consumer.
$ws.Rows(15).Insert()
$ws.Range("B15").Value = "CRM_client_<paris>"
$ws.Range("D15").Value = "BO"

[void]$ws.Range("C6").Select()
